# Updated cryptos list on Wed Apr 19 04:44:17 UTC 2023 with GitHub Actions
# Refresh Price (D) / Volume(1h) (E) figures for each coin row, and shift
# BabyDogeCoin into row 49 (pushing ThetaToken/EOS down one row, dropping
# Cronos from the bottom of the table).
# Leading "'" on each literal forces Excel to store the value as text,
# matching the original inline-string cells (prices like "341.57" would
# otherwise be auto-coerced to numbers).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.293.85"
$ws.Range("E2").Value = "'  +1.97%  "

$ws.Range("D3").Value = "'2.094.00"
$ws.Range("E3").Value = "'  -0.23%  "

$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "'  -0.73%  "

$ws.Range("D5").Value = "'341.57"
$ws.Range("E5").Value = "'  -0.61%  "

$ws.Range("E6").Value = "'  -0.58%  "

$ws.Range("D7").Value = "'0.5308"
$ws.Range("E7").Value = "'  +2.25%  "

$ws.Range("D8").Value = "'0.4382"
$ws.Range("E8").Value = "'  -0.01%  "

$ws.Range("D9").Value = "'54.15"
$ws.Range("E9").Value = "'  +0.60%  "

$ws.Range("D10").Value = "'0.09358"
$ws.Range("E10").Value = "'  +1.38%  "

$ws.Range("D11").Value = "'1.175"
$ws.Range("E11").Value = "'  +0.67%  "

$ws.Range("D12").Value = "'24.68"
$ws.Range("E12").Value = "'  +0.16%  "

$ws.Range("D13").Value = "'8.581"
$ws.Range("E13").Value = "'  +5.12%  "

$ws.Range("D14").Value = "'6.878"
$ws.Range("E14").Value = "'  +1.22%  "

$ws.Range("D15").Value = "'2.018.36"
$ws.Range("E15").Value = "'  -1.82%  "

$ws.Range("D16").Value = "'101.50"
$ws.Range("E16").Value = "'  -1.58%  "

$ws.Range("D17").Value = "'0.00001156"
$ws.Range("E17").Value = "'  +0.23%  "

$ws.Range("E18").Value = "'  -0.54%  "

$ws.Range("D19").Value = "'21.11"
$ws.Range("E19").Value = "'  +0.54%  "

$ws.Range("D20").Value = "'0.06723"
$ws.Range("E20").Value = "'  +0.89%  "

$ws.Range("D21").Value = "'6.332"
$ws.Range("E21").Value = "'  +1.91%  "

$ws.Range("D22").Value = "'1.001"
$ws.Range("E22").Value = "'  -0.64%  "

$ws.Range("D23").Value = "'30.281.70"
$ws.Range("E23").Value = "'  +1.87%  "

$ws.Range("D24").Value = "'12.47"
$ws.Range("E24").Value = "'  -0.57%  "

$ws.Range("D25").Value = "'2.317"
$ws.Range("E25").Value = "'  +0.55%  "

$ws.Range("D26").Value = "'21.81"
$ws.Range("E26").Value = "'  -0.51%  "

$ws.Range("D27").Value = "'6.939"
$ws.Range("E27").Value = "'  +9.12%  "

$ws.Range("D28").Value = "'162.28"
$ws.Range("E28").Value = "'  +0.17%  "

$ws.Range("D29").Value = "'2.506"
$ws.Range("E29").Value = "'  +0.38%  "

$ws.Range("D30").Value = "'133.72"
$ws.Range("E30").Value = "'  +0.04%  "

$ws.Range("D31").Value = "'1.129"
$ws.Range("E31").Value = "'  +0.03%  "

$ws.Range("D32").Value = "'0.1054"
$ws.Range("E32").Value = "'  +0.19%  "

$ws.Range("D33").Value = "'1.667"
$ws.Range("E33").Value = "'  -2.05%  "

$ws.Range("D34").Value = "'6.247"
$ws.Range("E34").Value = "'  +0.82%  "

$ws.Range("D35").Value = "'3.911"
$ws.Range("E35").Value = "'  -1.03%  "

$ws.Range("D36").Value = "'10.05"

$ws.Range("D37").Value = "'0.02625"
$ws.Range("E37").Value = "'  +1.92%  "

$ws.Range("D38").Value = "'0.06764"
$ws.Range("E38").Value = "'  +0.65%  "

$ws.Range("D39").Value = "'12.58"
$ws.Range("E39").Value = "'  +0.56%  "

$ws.Range("D40").Value = "'0.6955"
$ws.Range("E40").Value = "'  -0.44%  "

$ws.Range("D41").Value = "'1.341"
$ws.Range("E41").Value = "'  +0.91%  "

$ws.Range("D42").Value = "'0.2215"
$ws.Range("E42").Value = "'  +0.04%  "

$ws.Range("D43").Value = "'0.6782"
$ws.Range("E43").Value = "'  -0.19%  "

$ws.Range("D44").Value = "'2.352"
$ws.Range("E44").Value = "'  +0.89%  "

$ws.Range("D45").Value = "'14.24"
$ws.Range("E45").Value = "'  -0.42%  "

$ws.Range("D46").Value = "'1.001"
$ws.Range("E46").Value = "'  -0.55%  "

$ws.Range("D47").Value = "'1.282"
$ws.Range("E47").Value = "'  +6.75%  "

$ws.Range("D48").Value = "'3.635"
$ws.Range("E48").Value = "'  +0.39%  "

$ws.Range("B49").Value = "'BabyDogeCoin"
$ws.Range("C49").Value = "'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "'0.00000000348"
$ws.Range("E49").Value = "'  -2.92%  "

$ws.Range("B50").Value = "'ThetaToken"
$ws.Range("C50").Value = "'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D50").Value = "'1.206"
$ws.Range("E50").Value = "'  +4.17%  "

$ws.Range("B51").Value = "'EOS"
$ws.Range("C51").Value = "'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D51").Value = "'1.212"
$ws.Range("E51").Value = "'  -0.49%  "
